$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: fill in "Camera Back" entry under High-Res (was blank + 0:00)
$ws.Range("A22").Value = "Camera Back"
$ws.Range("B22").Value = "Medium"
$ws.Range("C22").Value = "00:20"

# New rows 23-25: more High-Res parts
$ws.Range("A23").Value = "Hinge"
$ws.Range("B23").Value = "Small"
$ws.Range("C23").Value = "00:05"

$ws.Range("A24").Value = "Attatchment"
$ws.Range("B24").Value = "Medium"
$ws.Range("C24").Value = "00:22"

$ws.Range("A25").Value = "Camera Back"
$ws.Range("B25").Value = "Medium"
$ws.Range("C25").Value = "00:23"

# Row 27: new "UV Mapping" section header
$ws.Range("A27").Value = "UV Mapping"
$ws.Range("A27").Style = $ws.Range("A21").Style

# Row 28: UV Mapping entry
$ws.Range("A28").Value = "Camera"
$ws.Range("B28").Value = "Large"
$ws.Range("C28").Value = "02:00"

# Row 30: new "Texturing" section header
$ws.Range("A30").Value = "Texturing"
$ws.Range("A30").Style = $ws.Range("A21").Style

# Row 31: Texturing entry
$ws.Range("A31").Value = "Camera"
$ws.Range("B31").Value = "Large"
$ws.Range("C31").Value = "03:00"

# Match styles for new data rows to existing data rows (s=3 for A/B, s=4 for C)
$ws.Range("A23:B25").Style = $ws.Range("A19:B19").Style
$ws.Range("C23:C25").Style = $ws.Range("C19").Style
$ws.Range("A28:B28").Style = $ws.Range("A19:B19").Style
$ws.Range("C28").Style = $ws.Range("C19").Style
$ws.Range("A31:B31").Style = $ws.Range("A19:B19").Style
$ws.Range("C31").Style = $ws.Range("C19").Style

# Update selection / view to match final state
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("C31").Select()
